$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 186, shifting existing rows 186-193 down to 187-194
$ws.Rows.Item(186).Insert()

# Populate the new row 186 with the new weekly record
$ws.Range("A186").Value = 8
$ws.Range("B186").Value = "Terminal La Palmera de La Serena"
$ws.Range("C186").Value = "Coquimbo"
$ws.Range("D186").Value = 44939
$ws.Range("E186").Value = 4
$ws.Range("F186").Value = 100112001
$ws.Range("G186").Value = "Berenjena"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 480
$ws.Range("K186").Value = 11000
$ws.Range("L186").Value = 12000
$ws.Range("M186").Value = 11500
$ws.Range("N186").Value = "$/caja 40 unidades"
$ws.Range("O186").Value = "Región de Arica y Parinacota"
$ws.Range("P186").Value = 288
$ws.Range("Q186").Value = 40
$ws.Range("R186").Value = "Hortaliza"
